{"js": "// Remove the first occurrence of the \"[PUMP:TBD:1]\" / \"BOLUS:SRS:2\"\n// paragraph pair from the body, leaving the second (identical) pair\n// untouched \u2014 matches the author's \"commented out the error line\" edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet pumpIndex = -1;\nlet bolusIndex = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const text = (p.text || \"\").trim();\n\n  if (pumpIndex === -1 && text === \"[PUMP:TBD:1]\") {\n    pumpIndex = i;\n    continue;\n  }\n\n  if (pumpIndex !== -1 && bolusIndex === -1 && text === \"BOLUS:SRS:2\" && p.style === \"List Bullet\") {\n    bolusIndex = i;\n    break;\n  }\n}\n\nif (pumpIndex !== -1 && bolusIndex !== -1) {\n  // Delete bottom-up so indices stay valid.\n  paragraphs.items[bolusIndex].delete();\n  paragraphs.items[pumpIndex].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the first occurrence of the \"[PUMP:TBD:1]\" / \"BOLUS:SRS:2\"\n# paragraph pair from the document, leaving the second (identical) pair\n# untouched \u2014 matches the author's \"commented out the error line\" edit.\n\n$d = $word.ActiveDocument\n\n$pumpIndex = -1\n$bolusIndex = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n\n    if ($pumpIndex -eq -1 -and $t -eq \"[PUMP:TBD:1]\") {\n        $pumpIndex = $i\n        continue\n    }\n\n    if ($pumpIndex -ne -1 -and $bolusIndex -eq -1 -and $t -eq \"BOLUS:SRS:2\" -and $p.Style.NameLocal -eq \"List Bullet\") {\n        $bolusIndex = $i\n        break\n    }\n}\n\nif ($pumpIndex -ne -1 -and $bolusIndex -ne -1) {\n    # Delete bottom-up so indices stay valid.\n    $d.Paragraphs.Item($bolusIndex).Range.Delete()\n    $d.Paragraphs.Item($pumpIndex).Range.Delete()\n}\n"}
